# Add data organization files for LEGOv2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("samples_retained")

# New row of data for the LEGOv2 dataset
$ws.Range("A20").Value = "LEGOv2"
$ws.Range("B20").Value = "spon."
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = 797
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = "English"
$ws.Range("H20").Value = "from an automated bus info service"

# Re-fill the "n" total formula down through row 28 (the new LEGOv2 row plus
# a handful of trailing blank rows), same formula already used in G2:G19.
$ws.Range("G2:G28").Formula = "=IF(OR(ISBLANK(C2), ISBLANK(D2),ISBLANK(E2)), """", SUM(C2:E2))"

# Update the selection to match the authored state
$ws.Range("G21").Select()
